# New crime data collected — updates the weekly CompStat 34th Precinct sheet:
#   - report header: volume number 15 -> 16, week-of dates shift by one week
#   - rows 15-29 (Murder .. Hate Crimes, TOTAL, Petit/Misd Larceny, UCR Rape*, etc.):
#     refreshed weekly/28-day/YTD/2yr counts and their derived % changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------------
# Header (shared, rich-text strings): only the embedded numbers change, the
# surrounding run text/formatting ("Volume ", "30", "   Number  ", etc. and
# "Report Covering the Week  " / "  Through  ") stays as-is.
# ---------------------------------------------------------------------------
$volCell = $ws.Range("A8")                      # "Volume 30   Number  15"
$volCell.Characters(21, 2).Text = "16"

$weekCell = $ws.Range("C9")                     # "...Week  4/10/2023  Through  4/16/2023"
$weekCell.Characters(27, 9).Text = "4/17/2023"
$weekCell.Characters(47, 9).Text = "4/23/2023"

# ---------------------------------------------------------------------------
# Helper: some cells flip between a numeric value and the literal text
# placeholders used for "no data" ("0" / "***.*", shared strings backing
# style s=14). Setting Value2 to a leading-apostrophe string forces text
# (mirrors typing '0 into Excel) while a plain numeric assignment is used
# for normal counts/percentages. PasteSpecial(formats) afterwards re-applies
# the correct existing cell style (xf index) from a stable, untouched cell
# that already carries it, since a plain assignment can otherwise reset the
# number format.
# ---------------------------------------------------------------------------
function Set-TextPlaceholder($addr, $text, $styleSrc) {
    $ws.Range($addr).Value2 = "'" + $text
    $ws.Range($styleSrc).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($fmt) | Out-Null
}
function Set-NumericWithStyle($addr, $value, $styleSrc) {
    $ws.Range($addr).Value2 = $value
    $ws.Range($styleSrc).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($fmt) | Out-Null
}

# Stable reference cells for each style, never themselves overwritten below.
$styleTextRef = "A14"   # s=14 (text "0"/"***.*" placeholder style)
$styleCountRef = "G15"  # s=16 (bold numeric count style)
$stylePctRef = "L14"    # s=15 (numeric percentage style)

# ---------------------------------------------------------------------------
# Row 15 — Rape: weekly count/% now "no data" placeholders.
# ---------------------------------------------------------------------------
Set-TextPlaceholder "D15" "0"    $styleTextRef
Set-TextPlaceholder "E15" "***.*" $styleTextRef
$ws.Range("N15").Value2 = -89.285714285714

# ---------------------------------------------------------------------------
# Row 16 — Robbery
# ---------------------------------------------------------------------------
$ws.Range("D16").Value2 = 1
$ws.Range("E16").Value2 = 100
$ws.Range("F16").Value2 = 11
$ws.Range("G16").Value2 = 13
$ws.Range("H16").Value2 = -15.384615384615
$ws.Range("I16").Value2 = 70
$ws.Range("J16").Value2 = 70
$ws.Range("K16").Value2 = 0
$ws.Range("L16").Value2 = 52.173913043478
$ws.Range("M16").Value2 = -27.083333333333
$ws.Range("N16").Value2 = -85.138004246284

# ---------------------------------------------------------------------------
# Row 17 — Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value2 = 11
$ws.Range("D17").Value2 = 6
$ws.Range("E17").Value2 = 83.333333333333
$ws.Range("F17").Value2 = 26
$ws.Range("G17").Value2 = 33
$ws.Range("H17").Value2 = -21.212121212121
$ws.Range("I17").Value2 = 91
$ws.Range("J17").Value2 = 89
$ws.Range("K17").Value2 = 2.247191011235
$ws.Range("L17").Value2 = 15.189873417721
$ws.Range("M17").Value2 = 62.5
$ws.Range("N17").Value2 = -64.313725490196

# ---------------------------------------------------------------------------
# Row 18 — Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value2 = 4
$ws.Range("E18").Value2 = 100
$ws.Range("F18").Value2 = 13
$ws.Range("G18").Value2 = 12
$ws.Range("H18").Value2 = 8.333333333333
$ws.Range("I18").Value2 = 53
$ws.Range("J18").Value2 = 55
$ws.Range("K18").Value2 = -3.636363636363
$ws.Range("L18").Value2 = 39.473684210526
$ws.Range("M18").Value2 = -11.666666666666
$ws.Range("N18").Value2 = -91.368078175895

# ---------------------------------------------------------------------------
# Row 19 — Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value2 = 17
$ws.Range("D19").Value2 = 14
$ws.Range("E19").Value2 = 21.428571428571
$ws.Range("F19").Value2 = 44
$ws.Range("G19").Value2 = 46
$ws.Range("H19").Value2 = -4.347826086956
$ws.Range("I19").Value2 = 159
$ws.Range("J19").Value2 = 187
$ws.Range("K19").Value2 = -14.973262032085
$ws.Range("L19").Value2 = -7.558139534883
$ws.Range("M19").Value2 = 62.244897959183
$ws.Range("N19").Value2 = -54.310344827586

# ---------------------------------------------------------------------------
# Row 20 — G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value2 = 9
$ws.Range("D20").Value2 = 5
$ws.Range("E20").Value2 = 80
$ws.Range("F20").Value2 = 19
$ws.Range("G20").Value2 = 23
$ws.Range("H20").Value2 = -17.391304347826
$ws.Range("I20").Value2 = 80
$ws.Range("J20").Value2 = 119
$ws.Range("K20").Value2 = -32.773109243697
$ws.Range("L20").Value2 = 48.148148148148
$ws.Range("M20").Value2 = 233.333333333333
$ws.Range("N20").Value2 = -87.711213517665

# ---------------------------------------------------------------------------
# Row 21 — TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value2 = 43
$ws.Range("D21").Value2 = 28
$ws.Range("E21").Value2 = 53.571428571428
$ws.Range("F21").Value2 = 113
$ws.Range("H21").Value2 = -11.71875
$ws.Range("I21").Value2 = 456
$ws.Range("J21").Value2 = 529
$ws.Range("K21").Value2 = -13.799621928166
$ws.Range("L21").Value2 = 13.715710723192
$ws.Range("M21").Value2 = 32.944606413994
$ws.Range("N21").Value2 = -80.944421228583

# ---------------------------------------------------------------------------
# Row 22 — Transit (weekly 28-day comparison now a "no data" placeholder)
# ---------------------------------------------------------------------------
$ws.Range("C22").Value2 = 1
Set-TextPlaceholder "D22" "0"     $styleTextRef
Set-TextPlaceholder "E22" "***.*" $styleTextRef
$ws.Range("I22").Value2 = 14
$ws.Range("K22").Value2 = 7.692307692307
$ws.Range("L22").Value2 = 100
$ws.Range("M22").Value2 = 55.555555555555

# ---------------------------------------------------------------------------
# Row 23 — Housing (D/E flip from placeholders back to real numbers)
# ---------------------------------------------------------------------------
Set-NumericWithStyle "D23" 1 $styleCountRef
Set-NumericWithStyle "E23" 0 $stylePctRef
$ws.Range("F23").Value2 = 2
$ws.Range("G23").Value2 = 4
$ws.Range("H23").Value2 = -50
$ws.Range("I23").Value2 = 11
$ws.Range("J23").Value2 = 11
$ws.Range("L23").Value2 = 37.5
$ws.Range("M23").Value2 = 22.222222222222

# ---------------------------------------------------------------------------
# Row 24 — Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value2 = 21
$ws.Range("D24").Value2 = 20
$ws.Range("E24").Value2 = 5
$ws.Range("F24").Value2 = 98
$ws.Range("G24").Value2 = 79
$ws.Range("H24").Value2 = 24.050632911392
$ws.Range("I24").Value2 = 382
$ws.Range("J24").Value2 = 418
$ws.Range("K24").Value2 = -8.612440191387
$ws.Range("L24").Value2 = 75.229357798165
$ws.Range("M24").Value2 = 102.116402116402

# ---------------------------------------------------------------------------
# Row 25 — Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("D25").Value2 = 8
$ws.Range("E25").Value2 = -25
$ws.Range("G25").Value2 = 37
$ws.Range("H25").Value2 = -5.405405405405
$ws.Range("I25").Value2 = 146
$ws.Range("J25").Value2 = 139
$ws.Range("K25").Value2 = 5.035971223021
$ws.Range("L25").Value2 = 29.203539823008
$ws.Range("M25").Value2 = -8.176100628930

# ---------------------------------------------------------------------------
# Row 26 — UCR Rape* (weekly counts now "no data" placeholders)
# ---------------------------------------------------------------------------
Set-TextPlaceholder "C26" "0"     $styleTextRef
Set-TextPlaceholder "D26" "0"     $styleTextRef
Set-TextPlaceholder "E26" "***.*" $styleTextRef
$ws.Range("L26").Value2 = -50

# ---------------------------------------------------------------------------
# Row 27 — Other Sex Crimes (C flips from placeholder back to a real number)
# ---------------------------------------------------------------------------
Set-NumericWithStyle "C27" 2 $styleCountRef
$ws.Range("G27").Value2 = 2
$ws.Range("H27").Value2 = 100
$ws.Range("I27").Value2 = 18
$ws.Range("K27").Value2 = 5.882352941176
$ws.Range("L27").Value2 = -5.263157894736

# ---------------------------------------------------------------------------
# Rows 28-29 — Shooting Vic. / Shooting Inc. (only the 2-yr % column moves)
# ---------------------------------------------------------------------------
$ws.Range("L28").Value2 = -87.5
$ws.Range("N28").Value2 = -98.529411764705
$ws.Range("L29").Value2 = -85.714285714285
$ws.Range("N29").Value2 = -98.412698412698

$excel.CutCopyMode = 0
Write-Output "Applied weekly CompStat refresh (34th Precinct)."
